$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new review row (row 9) ---------------------------------------
$ws.Range("A9").Value = "com.singleton.strechy"
$ws.Range("B9").Value = "stretchy"
$ws.Range("C9").Value = "eligitel@gmail.com"
$ws.Range("D9").Value = "ronenchen27@gmail.com"
$ws.Range("E9").Value = "27/5/2019 15:59"
$ws.Range("F9").Value = "Wonderful game with wonderful graphics and a variety of cars. Free offline game. One of the best I’ve played."

# --- Hyperlinks for the new email / recovery cells -----------------------
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:eligitel@gmail.com", "", "", "eligitel@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:ronenchen27@gmail.com", "", "", "ronenchen27@gmail.com")

# --- Match formatting used by the rest of the table ---------------------
# Column A / F use the "Mangal" font style (same as rows above)
$ws.Range("A9").Font.Name = "Mangal"
$ws.Range("F9").Font.Name = "Mangal"

# Columns C / D use the centered "Calibri" style (email / recovery columns).
# Applied after the hyperlink so it overrides the default hyperlink look
# (blue/underlined) with the same style used by the other rows.
$ws.Range("C9").Font.Name = "Calibri"
$ws.Range("C9").Font.Size = 11
$ws.Range("C9").Font.Color = 0
$ws.Range("C9").Font.Underline = $false
$ws.Range("C9").HorizontalAlignment = -4108

$ws.Range("D9").Font.Name = "Calibri"
$ws.Range("D9").Font.Size = 11
$ws.Range("D9").Font.Color = 0
$ws.Range("D9").Font.Underline = $false
$ws.Range("D9").HorizontalAlignment = -4108

# --- Update the active selection to the last cell of the new row ---------
$ws.Range("F9").Select() | Out-Null
